# Applies the "Copied documentation from branch new_version" edit.
#
# Real text-level changes identified from the diff (the rest of the
# diff is purely run-splitting/merging that does not change the
# visible text and is left untouched):
#
#   1. "http-server \u2013proxy" -> "http-server -\u2013proxy"
#      (a literal hyphen is inserted right before the existing en dash)
#   2. "http://localhost:8000" -> "http://[backend host]:8000"
#      (only the first occurrence, inside the hyperlink run; the other
#      "http://localhost:8000/end/" further down must stay untouched)
#   3. "Manueles starten von JSON Dateien" -> "Manuelles starten von
#      JSON Dateien" (typo fix in the heading)

$d = $word.ActiveDocument
$enDash = [char]0x2013

# 1. Insert a literal "-" right before the en dash in "http-server –proxy"
#    (there is an earlier, unrelated "http-server –g" elsewhere in the
#    document that must stay untouched, hence matching on "–proxy" too).
#    The match starts exactly at the beginning of the Consolas-styled run,
#    so locate it with Find and then insert just the "-" character at that
#    point instead of rewriting the whole matched range — this keeps the
#    surrounding run formatting (Consolas / grey) intact.
$f1 = $d.Content
$f1.Find.Text = "http-server " + $enDash + "proxy"
$f1.Find.Forward = $true
$f1.Find.Wrap = 1
$found1 = $f1.Find.Execute()
if ($found1) {
  $dashPos = $f1.Start + ("http-server ").Length
  $d.Range($dashPos, $dashPos).InsertBefore("-")
}

# 2. Replace "localhost" with "[backend host]" inside the first
#    "http://localhost:8000" hyperlink (the later
#    "http://localhost:8000/end/" occurrence must stay untouched).
#    Matching only "localhost" (not the "http://" prefix) keeps the match
#    away from the start of the run, which preserves the run's rPr
#    (Hyperlink style, Consolas font, grey color, no underline).
$f2 = $d.Content
$found2 = $f2.Find.Execute("http://localhost:8000", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
if ($found2) {
  $hostStart = $f2.Start + ("http://").Length
  $hostEnd = $hostStart + ("localhost").Length
  $d.Range($hostStart, $hostEnd).Text = "[backend host]"
}

# 3. Fix the heading typo "Manueles" -> "Manuelles".
$r3 = $d.Content
$r3.Find.Execute("Manueles starten von JSON Dateien", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "Manuelles starten von JSON Dateien", 1) | Out-Null
